$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("500CTOUT160", "Java Box (160oz)", "2", "94.99", "189.98"),
    @("5000TOUT96", "Java Box (96oz)", "2", "81.99", "163.98"),
    @("24510GCC", 'Cake Circle - 10" (Gold)', "2", "54.99", "109.98"),
    @("245CCGR2518", "Cake Board - Full Sheet", "1", "56.99", "56.99"),
    @("245882WB", "Box Cake - 8x8x2.5 (window)", "1", "68.60", "68.60"),
    @("130TONG10BLK", 'Black Plastic Tongs - 10.5"', "1", "37.99", "37.99"),
    @("5004CAFE", "Cup - Espresso (4oz)", "1", "35.49", "35.49"),
    @("43306HCUPC300", "Container - Muffin (6 Pack)", "1", "68.49", "68.49")
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    # SKU and Item are non-numeric text, so a plain assignment already stores
    # them as text (no quote-prefix style fork needed).
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    # Quantity / Cost Per / Total Cost are digit-only strings in the source
    # data; prefix with an apostrophe so Excel stores them as text (matching
    # the workbook's existing inlineStr/text convention) instead of coercing
    # them to numbers.
    $ws.Cells.Item($row, 3).Value = "'" + $data[$i][2]
    $ws.Cells.Item($row, 4).Value = "'" + $data[$i][3]
    $ws.Cells.Item($row, 5).Value = "'" + $data[$i][4]
}
